# Update countries & provincias Spain
# - Swap the "Chequia"/"Marruecos" rows (row 34 becomes Marruecos, row 35 becomes Chequia)
# - Refresh the covid-19 counters for several countries
# - Bump the "Datos actualizados" timestamp in A1

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Timestamp banner -------------------------------------------------
$ws.Range("A1").Value = "Datos actualizados a 18 de Octubre de 2020 a las 21:43"

# --- Countries whose alphabetical ordering swapped --------------------
# Row 34 held "Chequia" and row 35 held "Marruecos"; after the update
# "Marruecos" sorts before "Chequia", so the two rows trade country
# names (and the data below is refreshed to the new source numbers).
$ws.Range("A34").Value = "Marruecos"
$ws.Range("A35").Value = "Chequia"

# --- Refreshed case numbers --------------------------------------------
# Row 4: Estados Unidos
$ws.Range("B4").Value = 8376509
$ws.Range("C4").Value = 33844
$ws.Range("D4").Value = 5448765
$ws.Range("E4").Value = 2703185
$ws.Range("F4").Value = 0
$ws.Range("G4").Value = 277
$ws.Range("H4").Value = 224559

# Row 21: Alemania
$ws.Range("B21").Value = 366904
$ws.Range("C21").Value = 5171
$ws.Range("D21").Value = 290000
$ws.Range("E21").Value = 67038
$ws.Range("F21").Value = 0
$ws.Range("G21").Value = 13
$ws.Range("H21").Value = 9866

# Row 34: Marruecos (new data)
$ws.Range("B34").Value = 173632
$ws.Range("C34").Value = 2721
$ws.Range("D34").Value = 143972
$ws.Range("E34").Value = 26732
$ws.Range("F34").Value = 0
$ws.Range("G34").Value = 50
$ws.Range("H34").Value = 2928

# Row 35: Chequia (new data)
$ws.Range("B35").Value = 171487
$ws.Range("C35").Value = 2660
$ws.Range("D35").Value = 69090
$ws.Range("E35").Value = 100995
$ws.Range("F35").Value = 0
$ws.Range("G35").Value = 50
$ws.Range("H35").Value = 1402

# Row 52
$ws.Range("B52").Value = 89137
$ws.Range("C52").Value = 703
$ws.Range("D52").Value = 42649
$ws.Range("E52").Value = 45136
$ws.Range("F52").Value = 0
$ws.Range("G52").Value = 6
$ws.Range("H52").Value = 1352

# Row 105
$ws.Range("B105").Value = 11518
$ws.Range("C105").Value = 40
$ws.Range("D105").Value = 10427
$ws.Range("E105").Value = 1021
$ws.Range("F105").Value = 0

# Row 140
$ws.Range("B140").Value = 4322
$ws.Range("C140").Value = 18
$ws.Range("D140").Value = 3983
$ws.Range("E140").Value = 305
$ws.Range("F140").Value = 0

# Row 189
$ws.Range("B189").Value = 265
$ws.Range("C189").Value = 5
$ws.Range("D189").Value = 217
$ws.Range("E189").Value = 46
$ws.Range("F189").Value = 0
